$wb = $excel.ActiveWorkbook

# ------------------------------------------------------------------
# 1. Insert a new "2022-Q1" sheet right before the "总计" sheet. The
#    cleanest way to get an exact structural/format clone (sheetPr,
#    styles, header row, dimension, …) of the other quarterly sheets
#    is to copy the "2021-Q4" sheet itself, then trim it down to a
#    single data row and overwrite that row's values.
# ------------------------------------------------------------------
$totalSheet = $wb.Worksheets.Item("总计")
$template   = $wb.Worksheets.Item("2021-Q4")

$template.Copy($totalSheet)

# Re-fetch handles by name/position: object references captured
# before the Copy() call can end up pointing at the wrong sheet once
# the sheet collection shifts (the COM wrapper resolves sheet handles
# positionally), so grab fresh ones now.
$totalSheet = $wb.Worksheets.Item("总计")
$q1 = $wb.Worksheets.Item("2021-Q4 (2)")
$q1.Name = "2022-Q1"

# The template had 6 data rows (rows 2-7); keep only one.
$q1.Rows("3:7").Delete()

# Overwrite the single data row with the 2022-Q1 fund holding.
# (B-G match the existing sheets' convention of storing the figures
# as text, so force a text number format before assigning them.)
$q1.Range("B2:G2").NumberFormat = "@"
$q1.Range("A2").Value = 0
$q1.Range("B2").Value = "007152"
$q1.Range("C2").Value = "诺德策略精选混合"
$q1.Range("D2").Value = "0.42"
$q1.Range("E2").Value = "93.32"
$q1.Range("F2").Value = "4.24"
$q1.Range("G2").Value = "0.0178"
$q1.Range("H2").Value = 8

# ------------------------------------------------------------------
# 2. Update the "总计" sheet: insert a new top data row for 2022-Q1
#    and shift the existing rows down, renumbering column A.
# ------------------------------------------------------------------
$old = $totalSheet.Range("A2:D3").Value()
$totalSheet.Range("A3:D4").Value = $old

$totalSheet.Range("A2").Copy()
$totalSheet.Range("A4").PasteSpecial(-4122)

$totalSheet.Range("A2").Value = 0
$totalSheet.Range("B2").Value = "2022-Q1"
$totalSheet.Range("C2").Value = 1
$totalSheet.Range("D2").Value = 0.02

$totalSheet.Range("A3").Value = 1
$totalSheet.Range("A4").Value = 2

Write-Output "done"
